$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the stray _GoBack bookmark that used to sit right after
#    "be calculated based on the user entered scale)"
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Merge runs that were split across multiple <w:r> elements into a
#    single run (the visible text doesn't change, only the run
#    structure - a simple self Find/Replace collapses the runs).
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "  -h                   : print list of option and exit",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "  -h                   : print list of option and exit", 2)

$d.Content.Find.Execute(
    "  -b                   : batch processing, save detection result and exit (do not run the annotator)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "  -b                   : batch processing, save detection result and exit (do not run the annotator)", 2)

$d.Content.Find.Execute(
    "  -s <directory>       : directory to store the result files (this enable more detailed results to be stored)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "  -s <directory>       : directory to store the result files (this enable more detailed results to be stored)", 2)

$d.Content.Find.Execute(
    " <directory>       : save intermediate images in subdirectory (used only for debugging)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " <directory>       : save intermediate images in subdirectory (used only for debugging)", 2)

# ------------------------------------------------------------------
# 3. Append the new "BUDAS database connector" section at the very
#    end of the document (before the final section break).
# ------------------------------------------------------------------
$paraCount = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($paraCount)

# blank paragraph
$lastPara.Range.InsertParagraphAfter()

# "BUDAS database connector" heading paragraph
$p2 = $d.Paragraphs($d.Paragraphs.Count)
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs($d.Paragraphs.Count)
$p3.Range.Text = "BUDAS database connector"

# Final descriptive paragraph: two runs (second one carries a
# lastRenderedPageBreak) followed by the relocated _GoBack bookmark.
$insertPoint = $d.Content
$insertPoint.Collapse(0)

$xml = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">We have written a small python script (loadAnnotator.py) that allow you to transfer the information extracted from the annotator to the BUDAS database. For this initial version, you will have to go into the </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">code to update the user/password information for the database and the location of the BUDAS_output.txt file. </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@

$insertPoint.InsertXML($xml) | Out-Null

Write-Host "Edit complete. Paragraph count:" $d.Paragraphs.Count
